# Update MicroUI diagrams in the VDG (replace usage of the 'platform' terminology)
#
# 1) Refresh the cached text of the auto date fields (footer "datetimeFigureOut"
#    fields) on the notes master, the slide master and every slide layout.
# 2) On slide 3, rename the two "Platform" labels to "VEE Port" and nudge/resize
#    their text boxes to match the new (shorter) label.

$p = $ppt.ActivePresentation

function Update-DatePlaceholder($shapes, [string]$newText) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $sh = $shapes.Item($i)
        if ($sh.Name -like "Date Placeholder*") {
            $sh.TextFrame.TextRange.Text = $newText
        }
    }
}

# --- 1a. Notes master date field: 23/07/07 -> 12/1/2025 ---
$nm = $p.NotesMaster
Update-DatePlaceholder $nm.Shapes "12/1/2025"

# --- 1b. Slide master date field: 07/07/2023 -> 01/12/2025 ---
$m = $p.SlideMaster
Update-DatePlaceholder $m.Shapes "01/12/2025"

# --- 1c. Every slide layout's date field: 07/07/2023 -> 01/12/2025 ---
for ($li = 1; $li -le $m.CustomLayouts.Count; $li++) {
    $cl = $m.CustomLayouts.Item($li)
    Update-DatePlaceholder $cl.Shapes "01/12/2025"
}

# --- 2. Slide 3: "Platform" -> "VEE Port" (+ resize/reposition) ---
$s3 = $p.Slides.Item(3)
for ($i = 1; $i -le $s3.Shapes.Count; $i++) {
    $sh = $s3.Shapes.Item($i)
    if ($sh.HasTextFrame -and $sh.TextFrame.TextRange.Text -eq "Platform") {
        if ([Math]::Round($sh.Left) -eq 534) {
            # first shape: off x 6776595 -> 6895575, ext cx 725070 -> 720198
            $sh.Left = 542.9586614173228
            $sh.Width = 56.708505637007875
        } else {
            # second shape: off x 9272145 -> 9310443, ext cx 725070 -> 720198
            $sh.Left = 733.1057480314961
            $sh.Width = 56.708505637007875
        }
        $sh.TextFrame.TextRange.Text = "VEE Port"
    }
}
